# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Price/percentage columns are stored as text in the sheet (prices such as
# "26.635.57" use dots as thousands separators, not decimal points), so for
# any replacement value that Excel would otherwise auto-parse as a number
# we briefly force the cell to Text format, assign the value, then restore
# the default "Normal" style so no stray style attribute is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.635.57'
$ws.Cells.Item(2, 5).Value = '  +0.03%  '

$ws.Cells.Item(3, 4).Value = '1.597.10'
$ws.Cells.Item(3, 5).Value = '  +0.58%  '

$ws.Cells.Item(4, 5).Value = '  +0.04%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '211.62'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +0.24%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.516'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +1.20%  '

$ws.Cells.Item(7, 5).Value = '  +0.05%  '

$ws.Cells.Item(8, 5).Value = '  +0.24%  '

$ws.Cells.Item(9, 5).Value = '  -0.44%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '19.52'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -0.33%  '

$ws.Cells.Item(11, 5).Value = '  +0.47%  '

$ws.Cells.Item(12, 4).Value = '1.820.12'
$ws.Cells.Item(12, 5).Value = '  +0.57%  '

$ws.Cells.Item(13, 4).Value = '1.574.42'
$ws.Cells.Item(13, 5).Value = '  -0.91%  '

$ws.Cells.Item(14, 5).Value = '  -0.04%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '64.53'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -0.11%  '

$ws.Cells.Item(17, 4).Value = '26.613.74'
$ws.Cells.Item(17, 5).Value = '  +0.04%  '

$ws.Cells.Item(18, 4).Value = '0.0₃0732'
$ws.Cells.Item(18, 5).Value = '  +0.45%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '208.60'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -0.14%  '

$ws.Cells.Item(20, 5).Value = '  +0.04%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '6.97'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +3.69%  '

$ws.Cells.Item(22, 5).Value = '  +0.39%  '

$ws.Cells.Item(23, 5).Value = '  -2.28%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '8.89'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +0.30%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '145.23'
$ws.Cells.Item(25, 4).Style = 'Normal'

$ws.Cells.Item(26, 5).Value = '  +0.05%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '7.14'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -1.33%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.115'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +0.50%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '15.26'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -0.33%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.0507'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -0.35%  '

$ws.Cells.Item(31, 5).Value = '  +0.96%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.23'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -0.18%  '

$ws.Cells.Item(33, 5).Value = '  -4.37%  '

$ws.Cells.Item(34, 5).Value = '  +0.61%  '

$ws.Cells.Item(35, 4).Value = '1.280.12'

$ws.Cells.Item(36, 5).Value = '  +0.80%  '

$ws.Cells.Item(37, 5).Value = '  +1.09%  '

$ws.Cells.Item(38, 5).Value = '  -0.52%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.844'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +1.93%  '

$ws.Cells.Item(40, 5).Value = '  +0.08%  '

$ws.Cells.Item(41, 5).Value = '  +1.81%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '64.43'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +2.60%  '

$ws.Cells.Item(43, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.786'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -0.44%  '

$ws.Cells.Item(44, 2).Value = 'MXToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '2.19'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +1.40%  '

$ws.Cells.Item(45, 4).Value = '1.733.36'
$ws.Cells.Item(45, 5).Value = '  +0.59%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.911'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +8.34%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '89.83'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +0.31%  '

$ws.Cells.Item(48, 5).Value = '  -0.30%  '

$ws.Cells.Item(49, 5).Value = '  +4.57%  '

$ws.Cells.Item(50, 5).Value = '  +0.42%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '7.47'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -1.16%  '
